# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table with the latest scraped snapshot. Values are plain text (some
# "prices" use "." as a thousands separator, e.g. "52.230.54", and the
# percentages keep their original leading/trailing padding spaces), so
# every assignment is forced to Text via a leading apostrophe to stop
# Excel from reinterpreting/renormalising them as numbers; the style is
# reset back to Normal right after so no stray "quote prefix" style is
# left behind on cells that originally had no explicit style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''52.230.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +5.41%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.794.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +5.86%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.03%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''116.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +4.13%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''342.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = '''0.555'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +5.75%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -0.01%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  +6.11%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''42.10'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +6.27%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.0870'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +7.17%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  +0.41%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  +2.44%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''7.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +1.66%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.232.30'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +5.92%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''2.779.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +5.06%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''0.889'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +4.07%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''52.090.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +5.21%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''3.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +10.33%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''13.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +1.14%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''6.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +4.15%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.0₃0987'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +4.14%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''279.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +3.84%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''70.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.86%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  +10.36%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''26.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +3.11%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -0.08%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''10.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +0.48%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  +1.19%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  +3.12%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''34.83'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.67%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''50.48'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +1.72%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  +4.83%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''0.0828'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +2.50%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''2.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +4.54%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -0.07%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''18.98'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -0.40%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''5.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +1.13%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''3.26'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +5.11%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''2.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +26.88%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.0374'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +12.34%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''23.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +3.06%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  +4.29%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''2.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +3.98%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''126.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -2.38%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''2.104.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +1.93%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''3.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +1.66%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  +3.50%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''5.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +7.08%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.915'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +22.28%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''9.02'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +1.72%  '
$ws.Range("E51").Style = "Normal"
